{"js": "// Replace the outdated \"contratos\" pre-condition sentence with the\n// updated \"naturezas financeiras\" wording, keeping the leading\n// \"Deve haver\" and the trailing \"[Caso de Uso 45]\" reference intact.\nconst oldText = \"Deve haver registros de contratos previamente cadastrados [Caso de Uso 45]\";\nconst newText = \"Deve haver registros de naturezas financeiras previamente cadastradas [Caso de Uso 45]\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole sentence with the revised wording in a single pass,\n  // so the run keeps the original \"Tahoma\"/size-20/black formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n} else {\n  // Fallback: in case only the tail of the sentence can be located\n  // (e.g. differing whitespace), still perform the equivalent edit.\n  const tailResults = context.document.body.search(\n    \"registros de contratos previamente cadastrados\",\n    { matchCase: true }\n  );\n  tailResults.load(\"items\");\n  await context.sync();\n  if (tailResults.items.length > 0) {\n    tailResults.items[0].insertText(\n      \"registros de naturezas financeiras previamente cadastradas\",\n      Word.InsertLocation.replace\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the \"PR\u00c9-CONDI\u00c7\u00c3O(OES)\" bullet that used to reference contracts\n# (\"contratos\") so it references the correct prerequisite: previously\n# registered financial natures (\"naturezas financeiras\").\n$d = $word.ActiveDocument\n\n$oldText = \"Deve haver registros de contratos previamente cadastrados [Caso de Uso 45]\"\n$newText = \"Deve haver registros de naturezas financeiras previamente cadastradas [Caso de Uso 45]\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
